# Update the simulated loading_percent results for the 380 kV case (Case_1_49).
# Only the numeric result cells (columns B,D,E,F,G,H,J,K,O across rows 2-25) change;
# columns A, C, I, L, M, N and the header row are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 7.803613304306367
$ws.Range("D2").Value = 10.09061812575778
$ws.Range("E2").Value = 14.04781692612014
$ws.Range("F2").Value = 28.62160357906475
$ws.Range("G2").Value = 27.29186565023517
$ws.Range("H2").Value = 14.01401836301485
$ws.Range("J2").Value = 9.766914877609358
$ws.Range("K2").Value = 11.71492723955321
$ws.Range("O2").Value = 21.11886810691572

# Row 3
$ws.Range("B3").Value = 7.731792848869524
$ws.Range("D3").Value = 10.01612301753512
$ws.Range("E3").Value = 13.97569652429709
$ws.Range("F3").Value = 28.70745047824891
$ws.Range("G3").Value = 27.43222402868155
$ws.Range("H3").Value = 14.07503462071574
$ws.Range("J3").Value = 9.773389529311673
$ws.Range("K3").Value = 11.16787162906861
$ws.Range("O3").Value = 21.22558096903264

# Row 4
$ws.Range("B4").Value = 7.689193955303539
$ws.Range("D4").Value = 9.97167454347869
$ws.Range("E4").Value = 13.93406428223581
$ws.Range("F4").Value = 28.76823854570038
$ws.Range("G4").Value = 27.52912404557059
$ws.Range("H4").Value = 14.11505724493908
$ws.Range("J4").Value = 9.778935462884398
$ws.Range("K4").Value = 10.81725978402929
$ws.Range("O4").Value = 21.29639637224053

# Row 5
$ws.Range("B5").Value = 7.672230892039905
$ws.Range("D5").Value = 9.953901617963762
$ws.Range("E5").Value = 13.91777892088562
$ws.Range("F5").Value = 28.79503456085274
$ws.Range("G5").Value = 27.57128973777646
$ws.Range("H5").Value = 14.13201012922327
$ws.Range("J5").Value = 9.78159073811166
$ws.Range("K5").Value = 10.67083118880966
$ws.Range("O5").Value = 21.32658228888397

# Row 6
$ws.Range("B6").Value = 7.669438677005587
$ws.Range("D6").Value = 9.950971421344882
$ws.Range("E6").Value = 13.91511620965626
$ws.Range("F6").Value = 28.79960607753071
$ws.Range("G6").Value = 27.57845255628213
$ws.Range("H6").Value = 14.13486399961865
$ws.Range("J6").Value = 9.782055525635384
$ws.Range("K6").Value = 10.64630717457344
$ws.Range("O6").Value = 21.33167476267713

# Row 7
$ws.Range("B7").Value = 7.688963555732582
$ws.Range("D7").Value = 9.971433454311144
$ws.Range("E7").Value = 13.93384188092232
$ws.Range("F7").Value = 28.768591738873
$ws.Range("G7").Value = 27.52968188640797
$ws.Range("H7").Value = 14.11528327258385
$ws.Range("J7").Value = 9.778969672064056
$ws.Range("K7").Value = 10.81529916137417
$ws.Range("O7").Value = 21.29679809671204

# Row 8
$ws.Range("B8").Value = 7.778550291001308
$ws.Range("D8").Value = 10.0646734657479
$ws.Range("E8").Value = 14.02240834548409
$ws.Range("F8").Value = 28.64952309986419
$ws.Range("G8").Value = 27.33802613630555
$ws.Range("H8").Value = 14.03452577834771
$ws.Range("J8").Value = 9.768821628974006
$ws.Range("K8").Value = 11.52943953129202
$ws.Range("O8").Value = 21.15456229520595

# Row 9
$ws.Range("B9").Value = 7.965241798907524
$ws.Range("D9").Value = 10.25706891109554
$ws.Range("E9").Value = 14.21646649995549
$ws.Range("F9").Value = 28.48040719633528
$ws.Range("G9").Value = 27.04799375714536
$ws.Range("H9").Value = 13.89645972231422
$ws.Range("J9").Value = 9.761363350070907
$ws.Range("K9").Value = 12.80765303279181
$ws.Range("O9").Value = 20.91777097701611

# Row 10
$ws.Range("B10").Value = 8.107898146131973
$ws.Range("D10").Value = 10.40328149624158
$ws.Range("E10").Value = 14.370541785766
$ws.Range("F10").Value = 28.39576048838029
$ws.Range("G10").Value = 26.88820185823177
$ws.Range("H10").Value = 13.80739660943841
$ws.Range("J10").Value = 9.763438127051256
$ws.Range("K10").Value = 13.66632892709133
$ws.Range("O10").Value = 20.76966752771629

# Row 11
$ws.Range("B11").Value = 8.173723116741556
$ws.Range("D11").Value = 10.47064092665868
$ws.Range("E11").Value = 14.44291954667934
$ws.Range("F11").Value = 28.36591626293108
$ws.Range("G11").Value = 26.82727810359367
$ws.Range("H11").Value = 13.76956669467111
$ws.Range("J11").Value = 9.766013487925166
$ws.Range("K11").Value = 14.03857487337123
$ws.Range("O11").Value = 20.70794767982401

# Row 12
$ws.Range("B12").Value = 8.198760562942729
$ws.Range("D12").Value = 10.4962518168137
$ws.Range("E12").Value = 14.47063728343798
$ws.Range("F12").Value = 28.35586485607028
$ws.Range("G12").Value = 26.80591490353419
$ws.Range("H12").Value = 13.75562778711488
$ws.Range("J12").Value = 9.767222340087104
$ws.Range("K12").Value = 14.17683339706825
$ws.Range("O12").Value = 20.6853922684567

# Row 13
$ws.Range("B13").Value = 8.193363727559868
$ws.Range("D13").Value = 10.49073175129627
$ws.Range("E13").Value = 14.46465430319749
$ws.Range("F13").Value = 28.35797395202797
$ws.Range("G13").Value = 26.81043968368989
$ws.Range("H13").Value = 13.75861258712494
$ws.Range("J13").Value = 9.766951619055046
$ws.Range("K13").Value = 14.14717802417483
$ws.Range("O13").Value = 20.6902136053159

# Row 14
$ws.Range("B14").Value = 8.175780868830492
$ws.Range("D14").Value = 10.47274598468244
$ws.Range("E14").Value = 14.44519380344394
$ws.Range("F14").Value = 28.36506425409635
$ws.Range("G14").Value = 26.8254862387815
$ws.Range("H14").Value = 13.76841218525258
$ws.Range("J14").Value = 9.766108264991789
$ws.Range("K14").Value = 14.05000396839982
$ws.Range("O14").Value = 20.70607564669801

# Row 15
$ws.Range("B15").Value = 8.165024635960151
$ws.Range("D15").Value = 10.46174209313527
$ws.Range("E15").Value = 14.43331345871978
$ws.Range("F15").Value = 28.36957016779295
$ws.Range("G15").Value = 26.83492546903946
$ws.Range("H15").Value = 13.77446506022868
$ws.Range("J15").Value = 9.765622078697302
$ws.Range("K15").Value = 13.99012841981788
$ws.Range("O15").Value = 20.71589805280997

# Row 16
$ws.Range("B16").Value = 8.103613094616328
$ws.Range("D16").Value = 10.39889486358088
$ws.Range("E16").Value = 14.36585614465164
$ws.Range("F16").Value = 28.39788556907735
$ws.Range("G16").Value = 26.89242146729744
$ws.Range("H16").Value = 13.80992294533998
$ws.Range("J16").Value = 9.763302569747642
$ws.Range("K16").Value = 13.6416266247001
$ws.Range("O16").Value = 20.7738151317138

# Row 17
$ws.Range("B17").Value = 8.066161065654015
$ws.Range("D17").Value = 10.36054389475085
$ws.Range("E17").Value = 14.3250458413061
$ws.Range("F17").Value = 28.41747801685904
$ws.Range("G17").Value = 26.93071822018828
$ws.Range("H17").Value = 13.83236320373452
$ws.Range("J17").Value = 9.762296834374196
$ws.Range("K17").Value = 13.42307976122626
$ws.Range("O17").Value = 20.81079587843396

# Row 18
$ws.Range("B18").Value = 8.044708537235183
$ws.Range("D18").Value = 10.33856622689877
$ws.Range("E18").Value = 14.30178939000284
$ws.Range("F18").Value = 28.42956216247499
$ws.Range("G18").Value = 26.95385214385107
$ws.Range("H18").Value = 13.84552298378247
$ws.Range("J18").Value = 9.7618719972082
$ws.Range("K18").Value = 13.29565119800494
$ws.Range("O18").Value = 20.8325981091631

# Row 19
$ws.Range("B19").Value = 8.03746104041422
$ws.Range("D19").Value = 10.33113940197382
$ws.Range("E19").Value = 14.29395292441214
$ws.Range("F19").Value = 28.43379349446639
$ws.Range("G19").Value = 26.96187447947921
$ws.Range("H19").Value = 13.85002205792683
$ws.Range("J19").Value = 9.761754569846657
$ws.Range("K19").Value = 13.25221165756736
$ws.Range("O19").Value = 20.84007121063339

# Row 20
$ws.Range("B20").Value = 8.070138863517187
$ws.Range("D20").Value = 10.36461819189848
$ws.Range("E20").Value = 14.32936789231466
$ws.Range("F20").Value = 28.41530797498643
$ws.Range("G20").Value = 26.9265268182318
$ws.Range("H20").Value = 13.82994824040711
$ws.Range("J20").Value = 9.762388002750711
$ws.Range("K20").Value = 13.44652358179299
$ws.Range("O20").Value = 20.80680414053784

# Row 21
$ws.Range("B21").Value = 8.180942550280253
$ws.Range("D21").Value = 10.47802618285424
$ws.Range("E21").Value = 14.45090157018152
$ws.Range("F21").Value = 28.36294770567375
$ws.Range("G21").Value = 26.82102024741883
$ws.Range("H21").Value = 13.76552331471953
$ws.Range("J21").Value = 9.766349646722073
$ws.Range("K21").Value = 14.0786201352673
$ws.Range("O21").Value = 20.70139439052302

# Row 22
$ws.Range("B22").Value = 8.253995437000306
$ws.Range("D22").Value = 10.55273860826129
$ws.Range("E22").Value = 14.53212764809217
$ws.Range("F22").Value = 28.33601447030923
$ws.Range("G22").Value = 26.76202354435456
$ws.Range("H22").Value = 13.72567064844302
$ws.Range("J22").Value = 9.770299969261929
$ws.Range("K22").Value = 14.47595599153717
$ws.Range("O22").Value = 20.63726422864131

# Row 23
$ws.Range("B23").Value = 8.214955011215586
$ws.Range("D23").Value = 10.51281494819099
$ws.Range("E23").Value = 14.48861782231316
$ws.Range("F23").Value = 28.34972119648464
$ws.Range("G23").Value = 26.79259521472386
$ws.Range("H23").Value = 13.74673453092743
$ws.Range("J23").Value = 9.768067422675683
$ws.Range("K23").Value = 14.26535133842589
$ws.Range("O23").Value = 20.67105484861189

# Row 24
$ws.Range("B24").Value = 8.068340252859157
$ws.Range("D24").Value = 10.36277598018735
$ws.Range("E24").Value = 14.32741325014812
$ws.Range("F24").Value = 28.41628649642524
$ws.Range("G24").Value = 26.92841827405054
$ws.Range("H24").Value = 13.83103923988897
$ws.Range("J24").Value = 9.76234630771266
$ws.Range("K24").Value = 13.435930179594
$ws.Range("O24").Value = 20.8086071187396

# Row 25
$ws.Range("B25").Value = 7.913682216320778
$ws.Range("D25").Value = 10.2040962725142
$ws.Range("E25").Value = 14.16187899771909
$ws.Range("F25").Value = 28.51922639455005
$ws.Range("G25").Value = 27.11716489402914
$ws.Range("H25").Value = 13.93163726846015
$ws.Range("J25").Value = 9.762051694334531
$ws.Range("K25").Value = 12.47562491175969
$ws.Range("O25").Value = 20.97729999502784
